# Generate Report for Handoff
# - Refreshes the existing handoff entry (new guid-named file, new xlf names, new timestamps)
# - Appends a new row for a newly generated handoff file (ffffea43...)

$wb = $excel.ActiveWorkbook

$oldGuid = "8fe27c90-b8da-483b-bda5-03555ed4961e"
$newGuid = "c2abef80-4619-423b-95de-f222ad418178"
$newGuid2 = "ffffea43af6e-8ada-40de-a3a3-920232339166"

$oldXlfHash = "2cc13e4545c5f9127bd297223bb87d3ef53d1c0e"
$newXlfHash = "fce20ec6040f79f4076f816524777a37d98f7028"

$dateOverview = "2016-09-06 11:02:15"
$dateZhCn = "2016-09-06 11:01:55"
$dateDeDe = "2016-09-06 11:02:15"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58b273a9604b4d9510ab72c62cb365eb80ec9ae4/e2e/"

# -----------------------------------------------------------------
# Sheet "Overview"
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

# Update existing row 2 values
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "$repoBase$newGuid.md", "", "", "e2e\$newGuid.md") | Out-Null
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("D2").Value = ""
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = $dateOverview

# Add new row 3 (also expands the table range/dimension)
$loOverview.ListRows.Add() | Out-Null
$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "$repoBase$newGuid2.md", "", "", "e2e\$newGuid2.md") | Out-Null
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = $dateOverview

# -----------------------------------------------------------------
# Sheet "zh-cn"
# -----------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

# Update existing row 2 values
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$repoBase$newGuid.md", "", "", "$newGuid.md") | Out-Null
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("D2").Value = "e2e"
$wsZh.Range("E2").Value = "ht"
$wsZh.Range("F2").Value = "'False"
$wsZh.Range("G2").Value = "$newGuid.$newXlfHash.zh-cn.xlf"
$wsZh.Range("H2").Value = $dateZhCn
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"
$wsZh.Range("L2").Value = ""
$wsZh.Range("M2").Value = "'True"
$wsZh.Range("N2").Value = ""
$wsZh.Range("O2").Value = "'False"
$wsZh.Range("P2").Value = ""

# Add new row 3
$loZh.ListRows.Add() | Out-Null
$wsZh.Range("A3").Value = "$newGuid2.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$repoBase$newGuid2.md", "", "", "$newGuid2.md") | Out-Null
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("G3").Value = "$newGuid.$newXlfHash.zh-cn.xlf"
$wsZh.Range("H3").Value = $dateZhCn
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "'False"
$wsZh.Range("P3").Value = ""

# -----------------------------------------------------------------
# Sheet "de-de"
# -----------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

# Update existing row 2 values
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$repoBase$newGuid.md", "", "", "$newGuid.md") | Out-Null
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("D2").Value = "e2e"
$wsDe.Range("E2").Value = "ht"
$wsDe.Range("F2").Value = "'False"
$wsDe.Range("G2").Value = "$newGuid.$newXlfHash.de-de.xlf"
$wsDe.Range("H2").Value = $dateDeDe
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"
$wsDe.Range("L2").Value = ""
$wsDe.Range("M2").Value = "'True"
$wsDe.Range("N2").Value = ""
$wsDe.Range("O2").Value = "'False"
$wsDe.Range("P2").Value = ""

# Add new row 3
$loDe.ListRows.Add() | Out-Null
$wsDe.Range("A3").Value = "$newGuid2.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$repoBase$newGuid2.md", "", "", "$newGuid2.md") | Out-Null
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("G3").Value = "$newGuid.$newXlfHash.de-de.xlf"
$wsDe.Range("H3").Value = $dateDeDe
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "'False"
$wsDe.Range("P3").Value = ""

Write-Host "Handoff report generated."
